$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Cluster"
$ws.Range("B1").Value = "Active Cases"

# Cluster names for rows 2 through 60 (column A)
$labels = @(
  "3364 Assisi Centre Aged Care Rosanna",
  "3622 Olivet Care Aged Care Services Ringwood",
  "3633 Lifeview Emerald Glades Aged Care Emerald",
  "3961 Heritage Care Water Gardens Aged Care Facility Sydenham",
  "4167 Royal Freemasons Centennial Lodge Wantirna South",
  "4282 Villa Maria Catholic Homes (VMCH) Wantirna At-Home Aged Care Stud Road Wantirna",
  "AG Industries Pty Ltd Factory Thomastown",
  "Aintree Primary School Aintree",
  "Antonine College Cedar Campus Coburg",
  "Australian Meat Group Abattoir Dandenong South",
  "Bacchus Marsh Childcare and Kindergarten Centre Bacchus Marsh",
  "Baden Powell College Tarneit",
  "Covenant College Bell Post Hill",
  "Darul Ulum College of Victoria Fawkner October",
  "Domestic Freight Vessel Siem Aquamarine Geelong",
  "Gilly's Early Learning Centre Balaclava",
  "Gladstone Park Secondary College 29 Oct Gladstone Park",
  "Guardian Childcare & Education Moorabbin",
  "Hamlyn Banks Primary School Hamlyn Heights",
  "Hazelwood North Primary School Hazelwood North",
  "Islamic College of Melbourne Tarneit Oct Nov",
  "Lyndhurst Primary School Lyndhurst",
  "Master Poultry Group West Footscray",
  "Melbourne Metropolitan Remand Centre Ravenhall",
  "Morwell Park Primary School Morwell",
  "Nido Early School Woodend",
  "Northern Bay College Goldsworthy 9-12 Campus Corio",
  "Northern Bay College Wexford Campus Corio",
  "Northern Health Northern Hospital Epping Emergency Department Tier 1B",
  "Northern Health The Northern Hospital Epping",
  "Oakleigh South Primary School Oakleigh South",
  "Our Lady's Catholic Primary School Wangaratta",
  "Rutherglen Motor Inn and Walkabout Motel Rutherglen",
  "Sirius College Ibrahim Dellal Campus Sunshine",
  "Sirius College Shepparton Campus Shepparton",
  "Smartie Pants Early Learning and Development Diamond Creek",
  "St Brendans Primary School Shepparton",
  "St Clare's Primary School Officer",
  "St Georges Road Primary School Shepparton",
  "St Joseph's School Quarry Hill",
  "St Josephs Catholic Primary School Warragul",
  "St Louis de Montfort's School Aspendale",
  "St Paul's Primary School Sunshine West",
  "St Vincents Hospital Emergency Department Melbourne",
  "Stockdale Road Primary School Traralgon",
  "Story House Early Learning Epping October",
  "Sunbury Primary School Sunbury",
  "TLC Forest Lodge Residential Aged Care Frankston North",
  "TUROSI PTY LTD Thomastown",
  "Templestowe Park Primary School Templestowe",
  "The Lake Primary School Cabarita",
  "Tucker Road Bentleigh Primary School Bentleigh",
  "Warragul Regional College Warragul",
  "Werribee Mercy Hospital Emergency Department",
  "Western Health Sunshine Hospital Emergency Department St Albans",
  "Wodonga Primary School Wodonga",
  "Wodonga South Primary School Wodonga",
  "Wyndham Christian College Wyndham Vale",
  "Yeshivah College St Kilda East"
)

# Active case counts for rows 2 through 60 (column B)
$counts = @(19,12,10,24,22,11,17,11,10,10,34,14,20,10,10,10,14,11,10,28,24,10,14,10,53,11,11,31,14,14,16,10,22,10,12,20,10,12,14,32,11,13,12,11,33,12,11,10,11,30,19,10,16,21,12,15,13,12,15)

for ($i = 0; $i -lt $labels.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 1).Value = $labels[$i]
  $ws.Cells.Item($row, 2).Value = $counts[$i]
}

Write-Output "Done"
